# change tracing strategy and save wallet labels
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(
    "2024-10-03",
    "2024-09-27",
    "2024-09-26",
    "2024-10-05",
    "2024-10-04",
    "2024-10-02",
    "2024-10-01",
    "2024-09-30",
    "2024-08-02",
    "2024-08-12",
    "2024-09-25",
    "2024-09-24",
    "2024-09-23",
    "2024-09-20",
    "2024-09-18",
    "2024-09-16",
    "2024-09-14",
    "2024-09-13",
    "2024-09-12",
    "2024-09-11",
    "2024-09-10",
    "2024-09-09",
    "2024-09-08",
    "2024-09-07",
    "2024-09-06",
    "2024-08-09",
    "2024-07-24",
    "2024-09-17",
    "2024-08-05",
    "2024-03-09"
)

$startRow = 9
$endRow = $startRow + $dates.Length - 1
$newRange = $ws.Range("A$startRow`:A$endRow")

# Pre-format as text so the date-like strings ("2024-10-03", ...) are kept
# as literal text instead of being parsed into date serial numbers.
$newRange.NumberFormat = "@"

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
}

# Restore the default (unstyled) look so the new rows match the existing
# A2:A8 cells, which carry no explicit style.
$newRange.Style = "Normal"
